$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 13).Value = 257.5  # M28
$ws.Cells.Item(28, 14).Value = -2050.8462  # N28
$ws.Cells.Item(28, 10).Value = 1080.8462  # J28
$ws.Cells.Item(28, 8).Value = 563.6667  # H28
$ws.Cells.Item(28, 12).Value = 1080.8462  # L28
$ws.Cells.Item(28, 9).Value = 227.5  # I28
$ws.Cells.Item(28, 11).Value = 227.5  # K28

$ws.Cells.Item(62, 8).Value = 1319.375  # H62
$ws.Cells.Item(62, 12).Value = 1700  # L62
$ws.Cells.Item(62, 9).Value = 1091  # I62
$ws.Cells.Item(62, 11).Value = 1091  # K62
$ws.Cells.Item(62, 13).Value = -467  # M62
$ws.Cells.Item(62, 14).Value = -2948  # N62
$ws.Cells.Item(62, 10).Value = 1700  # J62

$ws.Cells.Item(63, 14).ClearContents()  # N63: was -11248
$ws.Cells.Item(63, 10).Value = 0  # J63
$ws.Cells.Item(63, 12).Value = 0  # L63
$ws.Cells.Item(63, 8).Value = 0  # H63

$ws.Cells.Item(65, 9).Value = 1091  # I65
$ws.Cells.Item(65, 11).Value = 5455  # K65
$ws.Cells.Item(65, 13).Value = -2335  # M65
$ws.Cells.Item(65, 14).Value = -14740  # N65
$ws.Cells.Item(65, 10).Value = 1700  # J65
$ws.Cells.Item(65, 12).Value = 8500  # L65
$ws.Cells.Item(65, 8).Value = 1319.375  # H65

$ws.Cells.Item(66, 14).ClearContents()  # N66: was -36240
$ws.Cells.Item(66, 10).Value = 0  # J66
$ws.Cells.Item(66, 12).Value = 0  # L66
$ws.Cells.Item(66, 8).Value = 0  # H66

$ws.Cells.Item(113, 8).Value = 54426  # H113
$ws.Cells.Item(113, 12).Value = 1910  # L113
$ws.Cells.Item(113, 9).Value = 101690.4  # I113
$ws.Cells.Item(113, 11).Value = 101690.4  # K113
$ws.Cells.Item(113, 13).Value = -98436.39999999999  # M113
$ws.Cells.Item(113, 14).Value = -8418  # N113
$ws.Cells.Item(113, 10).Value = 1910  # J113

$ws.Cells.Item(141, 8).Value = 2225.1667  # H141
$ws.Cells.Item(141, 9).Value = 2108.8696  # I141
$ws.Cells.Item(141, 11).Value = 6326.6088  # K141
$ws.Cells.Item(141, 13).Value = -1146.6088  # M141
$ws.Cells.Item(141, 14).Value = -25060  # N141
$ws.Cells.Item(141, 10).Value = 4900  # J141
$ws.Cells.Item(141, 12).Value = 14700  # L141

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 14).Value = -6340.3335  # N45
$ws.Cells.Item(45, 10).Value = 5586.3335  # J45
$ws.Cells.Item(45, 12).Value = 5586.3335  # L45
$ws.Cells.Item(45, 8).Value = 204151.8  # H45

$ws.Cells.Item(61, 14).Value = -2711.6  # N61
$ws.Cells.Item(61, 10).Value = 2287.6  # J61
$ws.Cells.Item(61, 12).Value = 2287.6  # L61
$ws.Cells.Item(61, 8).Value = 2287.6  # H61

$ws.Cells.Item(97, 11).Value = 59530.47  # K97
$ws.Cells.Item(97, 13).Value = -59034.47  # M97
$ws.Cells.Item(97, 14).Value = -2747.5  # N97
$ws.Cells.Item(97, 10).Value = 1755.5  # J97
$ws.Cells.Item(97, 8).Value = 48525.715  # H97
$ws.Cells.Item(97, 12).Value = 1755.5  # L97
$ws.Cells.Item(97, 9).Value = 59530.47  # I97

$ws.Cells.Item(106, 14).Value = -47459.75  # N106
$ws.Cells.Item(106, 10).Value = 44935.75  # J106
$ws.Cells.Item(106, 12).Value = 44935.75  # L106
$ws.Cells.Item(106, 8).Value = 44935.75  # H106

$ws.Cells.Item(122, 9).Value = 6293.3335  # I122
$ws.Cells.Item(122, 11).Value = 18880.0005  # K122
$ws.Cells.Item(122, 13).Value = -16430.0005  # M122
$ws.Cells.Item(122, 14).Value = -28900  # N122
$ws.Cells.Item(122, 10).Value = 8000  # J122
$ws.Cells.Item(122, 12).Value = 24000  # L122
$ws.Cells.Item(122, 8).Value = 6720  # H122

$ws.Cells.Item(136, 14).Value = -11962.8  # N136
$ws.Cells.Item(136, 10).Value = 2287.6  # J136
$ws.Cells.Item(136, 12).Value = 6862.799999999999  # L136
$ws.Cells.Item(136, 8).Value = 2287.6  # H136

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 9).Value = 51669.5  # I20
$ws.Cells.Item(20, 11).Value = 51669.5  # K20
$ws.Cells.Item(20, 13).Value = -51422.5  # M20
$ws.Cells.Item(20, 14).Value = -1518.8889  # N20
$ws.Cells.Item(20, 10).Value = 1024.8889  # J20
$ws.Cells.Item(20, 12).Value = 1024.8889  # L20
$ws.Cells.Item(20, 8).Value = 35952.207  # H20

$ws.Cells.Item(80, 9).Value = 614.5  # I80
$ws.Cells.Item(80, 11).Value = 614.5  # K80
$ws.Cells.Item(80, 14).Value = -5091.066699999999  # N80
$ws.Cells.Item(80, 13).Value = 383.5  # M80
$ws.Cells.Item(80, 10).Value = 3095.0667  # J80
$ws.Cells.Item(80, 12).Value = 3095.0667  # L80
$ws.Cells.Item(80, 8).Value = 2386.3333  # H80

$ws.Cells.Item(83, 13).Value = 1919.5  # M83
$ws.Cells.Item(83, 14).Value = -25459.3335  # N83
$ws.Cells.Item(83, 10).Value = 3095.0667  # J83
$ws.Cells.Item(83, 8).Value = 2386.3333  # H83
$ws.Cells.Item(83, 12).Value = 15475.3335  # L83
$ws.Cells.Item(83, 9).Value = 614.5  # I83
$ws.Cells.Item(83, 11).Value = 3072.5  # K83

$ws.Cells.Item(94, 8).Value = 603.4483  # H94
$ws.Cells.Item(94, 9).Value = 490.52942  # I94
$ws.Cells.Item(94, 11).Value = 490.52942  # K94
$ws.Cells.Item(94, 14).Value = -1665.4167  # N94
$ws.Cells.Item(94, 10).Value = 763.4167  # J94
$ws.Cells.Item(94, 13).Value = -39.52942000000002  # M94
$ws.Cells.Item(94, 12).Value = 763.4167  # L94

$ws.Cells.Item(99, 12).Value = 1479.8  # L99
$ws.Cells.Item(99, 8).Value = 1403.4166  # H99
$ws.Cells.Item(99, 9).Value = 1348.8572  # I99
$ws.Cells.Item(99, 11).Value = 1348.8572  # K99
$ws.Cells.Item(99, 13).Value = 149.1428000000001  # M99
$ws.Cells.Item(99, 14).Value = -4475.8  # N99
$ws.Cells.Item(99, 10).Value = 1479.8  # J99

$ws.Cells.Item(134, 9).Value = 2596.8333  # I134
$ws.Cells.Item(134, 11).Value = 7790.499899999999  # K134
$ws.Cells.Item(134, 13).Value = -5255.499899999999  # M134
$ws.Cells.Item(134, 8).Value = 2571.5527  # H134

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 12).Value = 54395.168  # L31
$ws.Cells.Item(31, 8).Value = 23522.816  # H31
$ws.Cells.Item(31, 11).Value = 933.29266  # K31
$ws.Cells.Item(31, 9).Value = 933.29266  # I31
$ws.Cells.Item(31, 14).Value = -54985.168  # N31
$ws.Cells.Item(31, 13).Value = -638.29266  # M31
$ws.Cells.Item(31, 10).Value = 54395.168  # J31

$ws.Cells.Item(34, 12).Value = 54395.168  # L34
$ws.Cells.Item(34, 8).Value = 23522.816  # H34
$ws.Cells.Item(34, 9).Value = 933.29266  # I34
$ws.Cells.Item(34, 11).Value = 933.29266  # K34
$ws.Cells.Item(34, 13).Value = -731.29266  # M34
$ws.Cells.Item(34, 14).Value = -54799.168  # N34
$ws.Cells.Item(34, 10).Value = 54395.168  # J34

$ws.Cells.Item(86, 8).Value = 2715.923  # H86
$ws.Cells.Item(86, 9).Value = 1866.6666  # I86
$ws.Cells.Item(86, 11).Value = 1866.6666  # K86
$ws.Cells.Item(86, 14).Value = -5216.7  # N86
$ws.Cells.Item(86, 10).Value = 2970.7  # J86
$ws.Cells.Item(86, 13).Value = -743.6666  # M86
$ws.Cells.Item(86, 12).Value = 2970.7  # L86

$ws.Cells.Item(89, 14).Value = -26085.5  # N89
$ws.Cells.Item(89, 10).Value = 2970.7  # J89
$ws.Cells.Item(89, 12).Value = 14853.5  # L89
$ws.Cells.Item(89, 8).Value = 2715.923  # H89
$ws.Cells.Item(89, 9).Value = 1866.6666  # I89
$ws.Cells.Item(89, 11).Value = 9333.333000000001  # K89
$ws.Cells.Item(89, 13).Value = -3717.333000000001  # M89

$ws.Cells.Item(99, 12).Value = 15322.875  # L99
$ws.Cells.Item(99, 8).Value = 13106.3  # H99
$ws.Cells.Item(99, 9).Value = 4240  # I99
$ws.Cells.Item(99, 11).Value = 4240  # K99
$ws.Cells.Item(99, 13).Value = -2742  # M99
$ws.Cells.Item(99, 14).Value = -18318.875  # N99
$ws.Cells.Item(99, 10).Value = 15322.875  # J99

$ws.Cells.Item(105, 8).Value = 2106.7827  # H105
$ws.Cells.Item(105, 9).Value = 2125.4736  # I105
$ws.Cells.Item(105, 11).Value = 2125.4736  # K105
$ws.Cells.Item(105, 13).Value = -378.4735999999998  # M105

$ws.Cells.Item(126, 14).Value = -50908.625  # N126
$ws.Cells.Item(126, 10).Value = 15322.875  # J126
$ws.Cells.Item(126, 8).Value = 13106.3  # H126
$ws.Cells.Item(126, 12).Value = 45968.625  # L126
$ws.Cells.Item(126, 9).Value = 4240  # I126
$ws.Cells.Item(126, 11).Value = 12720  # K126
$ws.Cells.Item(126, 13).Value = -10250  # M126

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(47, 8).Value = 157.42857  # H47
$ws.Cells.Item(47, 9).Value = 123.09091  # I47
$ws.Cells.Item(47, 11).Value = 369.27273  # K47
$ws.Cells.Item(47, 13).Value = 61.72727000000003  # M47

$ws.Cells.Item(112, 8).Value = 73398.14  # H112
$ws.Cells.Item(112, 9).Value = 201014.8  # I112
$ws.Cells.Item(112, 11).Value = 603044.3999999999  # K112
$ws.Cells.Item(112, 13).Value = -601936.3999999999  # M112

$ws.Cells.Item(136, 14).Value = -21798.9999  # N136
$ws.Cells.Item(136, 10).Value = 3866.3333  # J136
$ws.Cells.Item(136, 12).Value = 11598.9999  # L136
$ws.Cells.Item(136, 8).Value = 2771.5  # H136

$ws.Cells.Item(138, 12).Value = 9225  # L138
$ws.Cells.Item(138, 8).Value = 3088.3333  # H138
$ws.Cells.Item(138, 9).Value = 3115  # I138
$ws.Cells.Item(138, 11).Value = 9345  # K138
$ws.Cells.Item(138, 14).Value = -19505  # N138
$ws.Cells.Item(138, 10).Value = 3075  # J138
$ws.Cells.Item(138, 13).Value = -4205  # M138

$ws.Cells.Item(139, 14).Value = -19122.5879  # N139
$ws.Cells.Item(139, 10).Value = 2947.5293  # J139
$ws.Cells.Item(139, 12).Value = 8842.5879  # L139
$ws.Cells.Item(139, 8).Value = 2379.5  # H139
$ws.Cells.Item(139, 9).Value = 1306.5555  # I139
$ws.Cells.Item(139, 11).Value = 3919.6665  # K139
$ws.Cells.Item(139, 13).Value = 1220.3335  # M139

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 9).Value = 122089.414  # I70
$ws.Cells.Item(70, 11).Value = 122089.414  # K70
$ws.Cells.Item(70, 13).Value = -121819.414  # M70
$ws.Cells.Item(70, 14).Value = -7000.933  # N70
$ws.Cells.Item(70, 10).Value = 6460.933  # J70
$ws.Cells.Item(70, 8).Value = 67888.56  # H70
$ws.Cells.Item(70, 12).Value = 6460.933  # L70

$ws.Cells.Item(73, 11).Value = 122089.414  # K73
$ws.Cells.Item(73, 14).Value = -8332.933000000001  # N73
$ws.Cells.Item(73, 10).Value = 6460.933  # J73
$ws.Cells.Item(73, 13).Value = -121153.414  # M73
$ws.Cells.Item(73, 12).Value = 6460.933  # L73
$ws.Cells.Item(73, 8).Value = 67888.56  # H73
$ws.Cells.Item(73, 9).Value = 122089.414  # I73

$ws.Cells.Item(97, 11).Value = 166670500  # K97
$ws.Cells.Item(97, 13).Value = -166670004  # M97
$ws.Cells.Item(97, 14).Value = -3667  # N97
$ws.Cells.Item(97, 10).Value = 2675  # J97
$ws.Cells.Item(97, 8).Value = 100003370  # H97
$ws.Cells.Item(97, 12).Value = 2675  # L97
$ws.Cells.Item(97, 9).Value = 166670500  # I97

$ws.Cells.Item(107, 12).Value = 963.3333  # L107
$ws.Cells.Item(107, 8).Value = 425.90475  # H107
$ws.Cells.Item(107, 9).Value = 336.33334  # I107
$ws.Cells.Item(107, 11).Value = 336.33334  # K107
$ws.Cells.Item(107, 13).Value = 1583.66666  # M107
$ws.Cells.Item(107, 14).Value = -4803.3333  # N107
$ws.Cells.Item(107, 10).Value = 963.3333  # J107

$ws.Cells.Item(122, 9).Value = 10000  # I122
$ws.Cells.Item(122, 11).Value = 30000  # K122
$ws.Cells.Item(122, 13).Value = -27550  # M122
$ws.Cells.Item(122, 14).ClearContents()  # N122: was -10900
$ws.Cells.Item(122, 10).Value = 0  # J122
$ws.Cells.Item(122, 12).Value = 0  # L122
$ws.Cells.Item(122, 8).Value = 10000  # H122

$ws.Cells.Item(123, 14).Value = -14221.667  # N123
$ws.Cells.Item(123, 10).Value = 9321.666999999999  # J123
$ws.Cells.Item(123, 8).Value = 9321.666999999999  # H123
$ws.Cells.Item(123, 12).Value = 9321.666999999999  # L123

$ws.Cells.Item(126, 14).Value = -14721326  # N126
$ws.Cells.Item(126, 10).Value = 4905462  # J126
$ws.Cells.Item(126, 8).Value = 2031319.8  # H126
$ws.Cells.Item(126, 12).Value = 14716386  # L126
$ws.Cells.Item(126, 9).Value = 2513.5881  # I126
$ws.Cells.Item(126, 11).Value = 7540.7643  # K126
$ws.Cells.Item(126, 13).Value = -5070.7643  # M126

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 14).Value = -2692.8  # N40
$ws.Cells.Item(40, 10).Value = 2420.8  # J40
$ws.Cells.Item(40, 8).Value = 85559  # H40
$ws.Cells.Item(40, 12).Value = 2420.8  # L40

$ws.Cells.Item(82, 14).Value = -6714  # N82
$ws.Cells.Item(82, 10).Value = 5992  # J82
$ws.Cells.Item(82, 13).Value = -3639  # M82
$ws.Cells.Item(82, 12).Value = 5992  # L82
$ws.Cells.Item(82, 8).Value = 5660  # H82
$ws.Cells.Item(82, 9).Value = 4000  # I82
$ws.Cells.Item(82, 11).Value = 4000  # K82

$ws.Cells.Item(85, 12).Value = 5992  # L85
$ws.Cells.Item(85, 9).Value = 4000  # I85
$ws.Cells.Item(85, 11).Value = 4000  # K85
$ws.Cells.Item(85, 13).Value = -2752  # M85
$ws.Cells.Item(85, 14).Value = -8488  # N85
$ws.Cells.Item(85, 10).Value = 5992  # J85
$ws.Cells.Item(85, 8).Value = 5660  # H85

$ws.Cells.Item(140, 14).Value = -83308.25  # N140
$ws.Cells.Item(140, 10).Value = 72948.25  # J140
$ws.Cells.Item(140, 12).Value = 72948.25  # L140
$ws.Cells.Item(140, 8).Value = 72948.25  # H140

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(93, 14).Value = -35325.334  # N93
$ws.Cells.Item(93, 10).Value = 30333.334  # J93
$ws.Cells.Item(93, 12).Value = 30333.334  # L93
$ws.Cells.Item(93, 8).Value = 30333.334  # H93

$ws.Cells.Item(122, 14).Value = -16906  # N122
$ws.Cells.Item(122, 10).Value = 4002  # J122
$ws.Cells.Item(122, 12).Value = 12006  # L122
$ws.Cells.Item(122, 8).Value = 3201  # H122
